$d = $word.ActiveDocument

function Set-ParagraphRuns {
    param(
        [int]$ParaIndex,
        [string[]]$Texts
    )

    $p = $d.Paragraphs($ParaIndex)

    $runsXml = ""
    foreach ($t in $Texts) {
        $escaped = $t.Replace("&", "&amp;").Replace("<", "&lt;").Replace(">", "&gt;")
        $runsXml += '<w:r><w:t xml:space="preserve">' + $escaped + '</w:t></w:r>'
    }

    $bodyXml = '<w:p><w:pPr><w:pStyle w:val="FirstParagraph"/></w:pPr>' + $runsXml + '</w:p>'

    $pkg = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
           '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
           '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
           '<w:body>' + $bodyXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

    $p.Range.InsertXML($pkg) | Out-Null
}

# "Salut, les amis. Comment allez-vous ?" -> split into three runs
Set-ParagraphRuns 2 @("Salut, les amis.", " ", "Comment allez-vous ?")

# "J'ai 3 amis. J'ai 3 ans." -> split into five runs, appending " pas 2."
Set-ParagraphRuns 5 @("J’ai 3 amis.", " ", "J’ai 3 ans.", " ", "pas 2.")
